# Daily attendance processing - 2026-01-14 19:17:23
#
# In the "Recorded By" column (column G) of the session-analysis sheet,
# several cells list the recorder names as "System, dnasr281@gmail.com".
# Swap the order to "dnasr281@gmail.com, System" wherever that exact
# combination occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $used.Row + $used.Rows.Count - 1

# Bulk-read column G (Recorded By) once so we don't touch/materialize
# unrelated empty cells one-by-one.
$colRange = $ws.Range("G" + $firstRow + ":G" + $lastRow)
$vals = $colRange.Value2

$rowCount = $vals.GetLength(0)
for ($i = 1; $i -le $rowCount; $i++) {
    $cellValue = $vals[$i, 1]
    if ($cellValue -eq $oldValue) {
        $actualRow = $firstRow + $i - 1
        $ws.Cells.Item($actualRow, 7).Value = $newValue
    }
}
